# Updated cryptos list on Sun Sep 24 17:07:47 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.720.90"
$ws.Range("E2").Value = "  +0.05%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.598.93"
$ws.Range("E3").Value = "  -0.03%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.14%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.66"
$ws.Range("E5").Value = "  +0.10%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.09%  "

# Row 8 - Dogecoin
$ws.Range("E8").Value = "  -0.12%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.06%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.59"
$ws.Range("E10").Value = "  +0.45%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.82%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.823.77"
$ws.Range("E12").Value = "  +0.03%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.598.41"
$ws.Range("E13").Value = "  +0.06%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.55%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.24%  "

# Row 16 - Litecoin
$ws.Range("E16").Value = "  -0.30%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.0₃0740"
$ws.Range("E17").Value = "  -3.22%  "

# Row 18 - Dai
$ws.Range("E18").Value = "  +0.13%  "

# Row 19 - BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.35"
$ws.Range("E19").Value = "  -0.60%  "

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.15"
$ws.Range("E20").Value = "  +0.85%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +0.33%  "

# Row 22 - Toncoin
$ws.Range("E22").Value = "  -3.92%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  +0.65%  "

# Row 24 - Monero
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.71"
$ws.Range("E24").Value = "  +0.43%  "

# Row 25 - BinanceUSD
$ws.Range("E25").Value = "  -0.02%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  -0.31%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.63%  "

# Row 28 - EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.32"
$ws.Range("E28").Value = "  -0.04%  "

# Row 29 - Hedera
$ws.Range("E29").Value = "  -2.47%  "

# Row 30 - PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.16"
$ws.Range("E30").Value = "  +0.16%  "

# Row 31 - Filecoin
$ws.Range("E31").Value = "  +0.41%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +0.35%  "

# Row 33 - WEMIXToken
$ws.Range("E33").Value = "  +17.98%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.275.31"
$ws.Range("E34").Value = "  -0.69%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +1.69%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  +0.27%  "

# Row 37 - ImmutableX
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.591"
$ws.Range("E37").Value = "  -4.40%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -1.30%  "

# Row 39 - ARBITRUM
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.823"
$ws.Range("E39").Value = "  -0.47%  "

# Row 40 - FraxShare
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.47"
$ws.Range("E40").Value = "  +0.40%  "

# Row 41 - MXToken
$ws.Range("E41").Value = "  +0.24%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.775"
$ws.Range("E42").Value = "  -1.19%  "

# Row 43 - Aave
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.57"
$ws.Range("E43").Value = "  -1.09%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.735.30"
$ws.Range("E44").Value = "  +0.07%  "

# Row 45 - Quant
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.36"
$ws.Range("E45").Value = "  -0.79%  "

# Row 46 - RenderToken
$ws.Range("E46").Value = "  +0.49%  "

# Row 47 - Algorand
$ws.Range("E47").Value = "  +1.51%  "

# Row 48 - Cronos
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0513"
$ws.Range("E48").Value = "  +1.04%  "

# Row 49 - EnergySwap
$ws.Range("E49").Value = "  +2.79%  "

# Row 50 - USDD
$ws.Range("E50").Value = "  +0.07%  "

# Row 51 - Mantle
$ws.Range("E51").Value = "  +1.59%  "
